# Applies:
#  1) The table on slide 6 switches to table style {83B50A46-3AC8-4FF9-A17B-54FD31231B14}
#     (previously {D8BF2E58-6E85-444C-B3E0-5FB523C2E5A6}).
#  2) The deck's theme colour scheme changes from the "Integral" palette to the
#     classic "Office" palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 -------------------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{83B50A46-3AC8-4FF9-A17B-54FD31231B14}")
    }
}

# --- 2. Theme colour scheme: Integral -> Office --------------------------------
function HexToOleRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Order matches MsoThemeColorSchemeIndex 1..12:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$colorScheme = $p.Slides.Item(1).Master.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Length; $i++) {
    $colorScheme.Item($i).RGB = HexToOleRgb($officeColors[$i - 1])
}
